$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.179.15'
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.826.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.84%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9987'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6146'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.31%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07112'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.37%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2820'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.66'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.01%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07666'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.82%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.821.11'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.24%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.826'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.79%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.00001011'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6336'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.069.50'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '79.08'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.96%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.876'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.88%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.148.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '228.01'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.023'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.87'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.01%  '

$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1324'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.80%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.055'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.63'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.488'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.59%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06369'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.90%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.452'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.95%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.828'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.799'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.73%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.127'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.02%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.750'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6508'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.543'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.49%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.751'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.25%  '

$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.218.78'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.33%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.595'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.10%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01745'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9283'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9991'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.54%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.974.20'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.12%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.09'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.35%  '

$ws.Range("E47").Value = '  -0.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.628'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.53%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.623'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4560'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.66%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05517'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.72%  '
